# Localize the "empty" template (slide master, slide layouts, and the
# single title-slide placeholder slide) from en-US boilerplate to de-DE.
#
# For every shape we only touch:
#   - the shape's Name (p:cNvPr/@name)
#   - the cSld/@name of each slide layout
#   - the literal run text of the "Click to edit..." / "Second level" etc.
#     placeholder prompts (paragraph-by-paragraph, so pPr/lvl and
#     endParaRPr stay untouched)
# Field runs (date, slide number) and shape geometry are left alone.

$p = $ppt.ActivePresentation

# ---- text translations (exact paragraph text match) ----------------------
$textMap = @{
    "Click to edit Master title style"    = "Titelmasterformat durch Klicken bearbeiten"
    "Click to edit Master subtitle style" = "Formatvorlage des Untertitelmasters durch Klicken bearbeiten"
    "Click to edit Master text styles"    = "Textmasterformat bearbeiten"
    "Second level"                        = "Zweite Ebene"
    "Third level"                         = "Dritte Ebene"
    "Fourth level"                        = "Vierte Ebene"
    "Fifth level"                         = "Fünfte Ebene"
}

# ---- shape-name translations ----------------------------------------------
$nameMap = @{
    "Title 1"                          = "Titel 1"
    "Subtitle 2"                       = "Untertitel 2"
    "Title Placeholder 1"              = "Titelplatzhalter 1"
    "Vertical Title 1"                 = "Vertikaler Titel 1"
    "Content Placeholder 2"            = "Inhaltsplatzhalter 2"
    "Content Placeholder 3"            = "Inhaltsplatzhalter 3"
    "Content Placeholder 5"            = "Inhaltsplatzhalter 5"
    "Text Placeholder 2"               = "Textplatzhalter 2"
    "Text Placeholder 3"               = "Textplatzhalter 3"
    "Text Placeholder 4"               = "Textplatzhalter 4"
    "Vertical Text Placeholder 2"      = "Vertikaler Textplatzhalter 2"
    "Picture Placeholder 2"            = "Bildplatzhalter 2"
    "Date Placeholder 1"               = "Datumsplatzhalter 1"
    "Date Placeholder 2"               = "Datumsplatzhalter 2"
    "Date Placeholder 3"               = "Datumsplatzhalter 3"
    "Date Placeholder 4"               = "Datumsplatzhalter 4"
    "Date Placeholder 6"               = "Datumsplatzhalter 6"
    "Footer Placeholder 2"             = "Fußzeilenplatzhalter 2"
    "Footer Placeholder 3"             = "Fußzeilenplatzhalter 3"
    "Footer Placeholder 4"             = "Fußzeilenplatzhalter 4"
    "Footer Placeholder 5"             = "Fußzeilenplatzhalter 5"
    "Footer Placeholder 7"             = "Fußzeilenplatzhalter 7"
    "Slide Number Placeholder 3"       = "Foliennummernplatzhalter 3"
    "Slide Number Placeholder 4"       = "Foliennummernplatzhalter 4"
    "Slide Number Placeholder 5"       = "Foliennummernplatzhalter 5"
    "Slide Number Placeholder 6"       = "Foliennummernplatzhalter 6"
    "Slide Number Placeholder 8"       = "Foliennummernplatzhalter 8"
}

# ---- slide-layout (cSld) name translations --------------------------------
$layoutNameMap = @{
    "Title Slide"               = "Titelfolie"
    "Title and Vertical Text"   = "Titel und vertikaler Text"
    "Vertical Title and Text"   = "Vertikaler Titel und Text"
    "Title and Content"         = "Titel und Inhalt"
    "Section Header"            = "Abschnitts-überschrift"
    "Two Content"               = "Zwei Inhalte"
    "Comparison"                = "Vergleich"
    "Title Only"                = "Nur Titel"
    "Blank"                     = "Leer"
    "Content with Caption"      = "Inhalt mit Überschrift"
    "Picture with Caption"      = "Bild mit Überschrift"
}

function Rename-Shape($shape) {
    if ($nameMap.ContainsKey($shape.Name)) {
        $shape.Name = $nameMap[$shape.Name]
    }
}

function Translate-TextFrame($shape) {
    if (-not $shape.HasTextFrame) {
        return
    }
    $tr = $shape.TextFrame.TextRange
    $paras = $tr.Paragraphs()
    for ($j = 1; $j -le $paras.Count; $j++) {
        $para = $tr.Paragraphs($j, 1)
        $t = $para.Text
        if ($textMap.ContainsKey($t)) {
            $para.Text = $textMap[$t]
        }
    }
}

# ---- slide master ----------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    Rename-Shape $shape
}

# ---- slide layouts ----------------------------------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)

    if ($layoutNameMap.ContainsKey($layout.Name)) {
        $layout.Name = $layoutNameMap[$layout.Name]
    }

    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        Rename-Shape $shape
        Translate-TextFrame $shape
    }
}

# ---- slide 1 (title-slide placeholders, names only; text stays empty) -----
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    Rename-Shape $shape
}
